$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "38.777.11"
Set-TextValue 2 5 "  +1.71%  "
Set-TextValue 3 4 "2.098.29"
Set-TextValue 3 5 "  +0.32%  "
Set-TextValue 4 5 "  -0.05%  "
Set-TextValue 5 4 "229.00"
Set-TextValue 5 5 "  +0.04%  "
Set-TextValue 6 4 "0.618"
Set-TextValue 6 5 "  +0.89%  "
Set-TextValue 7 4 "61.85"
Set-TextValue 7 5 "  +2.20%  "
Set-TextValue 9 5 "  +2.24%  "
Set-TextValue 10 4 "0.0847"
Set-TextValue 10 5 "  +0.57%  "
Set-TextValue 11 5 "  +0.08%  "
Set-TextValue 12 4 "15.45"
Set-TextValue 12 5 "  +5.84%  "
Set-TextValue 13 4 "2.409.85"
Set-TextValue 13 5 "  +0.45%  "
Set-TextValue 14 4 "22.11"
Set-TextValue 14 5 "  -0.17%  "
Set-TextValue 15 5 "  +4.64%  "
Set-TextValue 16 4 "5.50"
Set-TextValue 16 5 "  +0.56%  "
Set-TextValue 17 4 "2.106.14"
Set-TextValue 17 5 "  +1.46%  "
Set-TextValue 18 4 "38.814.16"
Set-TextValue 18 5 "  +1.97%  "
Set-TextValue 19 4 "71.97"
Set-TextValue 19 5 "  +2.67%  "
Set-TextValue 20 4 "6.08"
Set-TextValue 20 5 "  +1.18%  "
$tmpVal1 = "0.0" + ([string][char]0x2083) + "0840"
Set-TextValue 21 4 $tmpVal1
Set-TextValue 21 5 "  +0.79%  "
Set-TextValue 22 4 "227.91"
Set-TextValue 22 5 "  +1.77%  "
Set-TextValue 23 4 "1.00"
Set-TextValue 23 5 "  +0.03%  "
Set-TextValue 24 5 "  -2.15%  "
Set-TextValue 25 5 "  +0.57%  "
Set-TextValue 26 4 "171.72"
Set-TextValue 26 5 "  +1.04%  "
Set-TextValue 27 4 "9.53"
Set-TextValue 27 5 "  +1.36%  "
Set-TextValue 28 5 "  +5.93%  "
Set-TextValue 29 5 "  +5.33%  "
Set-TextValue 30 4 "19.35"
Set-TextValue 30 5 "  +2.16%  "
Set-TextValue 31 4 "2.47"
Set-TextValue 31 5 "  +3.59%  "
Set-TextValue 32 5 "  +1.21%  "
Set-TextValue 33 4 "4.53"
Set-TextValue 33 5 "  +2.50%  "
Set-TextValue 34 4 "4.77"
Set-TextValue 34 5 "  +1.67%  "
Set-TextValue 35 4 "0.0621"
Set-TextValue 35 5 "  +2.57%  "
Set-TextValue 36 4 "6.51"
Set-TextValue 36 5 "  +1.43%  "
Set-TextValue 37 4 "2.40"
Set-TextValue 37 5 "  +0.00%  "
Set-TextValue 38 4 "3.58"
Set-TextValue 38 5 "  +1.60%  "
Set-TextValue 40 4 "18.13"
Set-TextValue 40 5 "  +0.46%  "
Set-TextValue 41 4 "0.0228"
Set-TextValue 41 5 "  +4.53%  "
Set-TextValue 42 4 "101.64"
Set-TextValue 42 5 "  +1.60%  "
Set-TextValue 43 4 "1.533.00"
Set-TextValue 43 5 "  -1.62%  "
Set-TextValue 44 5 "  -0.82%  "
Set-TextValue 45 4 "7.75"
Set-TextValue 45 5 "  +4.12%  "
Set-TextValue 46 5 "  -0.36%  "
Set-TextValue 47 5 "  +2.03%  "
Set-TextValue 48 4 "4.10"
Set-TextValue 48 5 "  -1.29%  "
Set-TextValue 49 5 "  +2.03%  "
Set-TextValue 50 5 "  -0.61%  "
Set-TextValue 51 4 "2.293.55"
Set-TextValue 51 5 "  +0.28%  "
